# Applies the two content changes described by the diff:
#   1. Move the "_GoBack" bookmark from the last paragraph ("equals() in
#      Recipe.java ...") to the end of the paragraph that contains
#      `newRecipe.setName("");` (right after that run).
#   2. Insert a new paragraph "Defects found from Inventory Tests" right
#      after the paragraph containing `recipeArray[recipeToEdit] = newRecipe;`.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Relocate the "_GoBack" bookmark.
# ---------------------------------------------------------------------
# Find the paragraph that starts with `newRecipe.setName("");`
# (NOTE: `-like` treats `[`/`]` as wildcard character-class syntax, so a
# plain `.StartsWith()` call is used instead to match the literal text.)
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text.StartsWith('newRecipe.setName("");')) {
        $targetPara = $para
        break
    }
}

if ($targetPara -eq $null) {
    throw "Could not locate the 'newRecipe.setName(`"`");' paragraph"
}

# The position right after the run's text (before the paragraph mark).
$endPos = $targetPara.Range.End - 1

# A genuinely zero-length Range right at that boundary is mishandled when
# passed straight to Bookmarks.Add, so nudge around it: insert a throw-away
# character there, wrap a bookmark around that single character, then strip
# the character back out again. The net effect is an empty "_GoBack"
# bookmark sitting exactly after the run, which also replaces (moves) any
# pre-existing "_GoBack" bookmark elsewhere in the document.
$insertionPoint = $d.Range($endPos, $endPos)
$insertionPoint.InsertAfter("X")
$markerRange = $d.Range($endPos, $endPos + 1)
$d.Bookmarks.Add("_GoBack", $markerRange)
$markerRange.Text = ""
Write-Output "Moved _GoBack bookmark after 'newRecipe.setName(`"`");'"

# ---------------------------------------------------------------------
# 2. Insert the new "Defects found from Inventory Tests" paragraph.
# ---------------------------------------------------------------------
$afterPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text.StartsWith('recipeArray[recipeToEdit] = newRecipe;')) {
        $afterPara = $para
        break
    }
}

if ($afterPara -eq $null) {
    throw "Could not locate the 'recipeArray[recipeToEdit] = newRecipe;' paragraph"
}

$afterPara.Range.InsertParagraphAfter()

$newParaIndex = $afterPara.Index + 1
$newPara = $d.Paragraphs($newParaIndex)
$newPara.Style = "Normal"
$newPara.Range.Text = "Defects found from Inventory Tests"
Write-Output "Inserted 'Defects found from Inventory Tests' paragraph"
